# Resource/excel/商城-商品配置.xlsx — "update resource by parse tool"
#
# Adjust the shop "money" prices for a few goods rows and tidy up a handful
# of stray, styled-but-empty header cells on row 1, then leave the selection
# on the last-touched cell (L12) as the tool would after finishing its pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (K) / Discount (L) columns: update the "money" JSON payloads ---
# Row 5 (goods id 10001, hero item): 100 -> 80
$ws.Range("K5").Value = '[{"money":"80"}]'
$ws.Range("L5").Value = '[{"money":"80"}]'

# Row 6 (goods id 11001, effect item): 100 -> 800
$ws.Range("K6").Value = '[{"money":"800"}]'
$ws.Range("L6").Value = '[{"money":"800"}]'

# Row 7 (goods id 12001, foot item): 100 -> 400
$ws.Range("K7").Value = '[{"money":"400"}]'
$ws.Range("L7").Value = '[{"money":"400"}]'

# --- Row 1: drop the leftover styled-but-empty cells in G/H/K/L ---
$ws.Range("G1").Clear()
$ws.Range("H1").Clear()
$ws.Range("K1").Clear()
$ws.Range("L1").Clear()

# --- Leave the selection where the parse tool left it ---
$ws.Range("L12").Select()
